$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cells whose new text looks like a plain number (e.g. "1.002") ---
# Pre-format as Text so Excel stores the literal string instead of
# silently converting it to a number, then restore the default style
# so no stray formatting is introduced.
$textCells = @("D4","D5","D7","D8","D9","D10","D11","D12","D14","D15","D16","D17","D18","D21","D24","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D49","D50","D51")
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D4").Value = '1.002'
$ws.Range("D5").Value = '312.50'
$ws.Range("D7").Value = '0.4879'
$ws.Range("D8").Value = '0.3796'
$ws.Range("D9").Value = '0.07325'
$ws.Range("D10").Value = '0.9135'
$ws.Range("D11").Value = '20.55'
$ws.Range("D12").Value = '0.07670'
$ws.Range("D14").Value = '5.470'
$ws.Range("D15").Value = '6.598'
$ws.Range("D16").Value = '91.01'
$ws.Range("D17").Value = '1.002'
$ws.Range("D18").Value = '0.000008773'
$ws.Range("D21").Value = '14.51'
$ws.Range("D24").Value = '10.73'
$ws.Range("D25").Value = '1.909'
$ws.Range("D26").Value = '153.72'
$ws.Range("D27").Value = '18.37'
$ws.Range("D28").Value = '2.146'
$ws.Range("D29").Value = '115.59'
$ws.Range("D30").Value = '4.892'
$ws.Range("D31").Value = '0.08915'
$ws.Range("D32").Value = '3.200'
$ws.Range("D33").Value = '1.220'
$ws.Range("D34").Value = '0.7668'
$ws.Range("D35").Value = '4.636'
$ws.Range("D37").Value = '2.528'
$ws.Range("D38").Value = '1.092'
$ws.Range("D39").Value = '0.05269'
$ws.Range("D40").Value = '0.5466'
$ws.Range("D41").Value = '2.977'
$ws.Range("D42").Value = '6.913'
$ws.Range("D43").Value = '8.514'
$ws.Range("D44").Value = '0.1517'
$ws.Range("D45").Value = '112.16'
$ws.Range("D46").Value = '10.60'
$ws.Range("D47").Value = '0.4788'
$ws.Range("D49").Value = '1.629'
$ws.Range("D50").Value = '67.37'
$ws.Range("D51").Value = '0.06046'

foreach ($ref in $textCells) {
    $ws.Range($ref).Style = "Normal"
}

# --- Cells whose new text is safe as-is (already non-numeric, e.g.
#     multi-dot price strings, percentages with padding spaces, or plain
#     text like coin names / URLs) ---
$ws.Range("D2").Value = '27.694.72'
$ws.Range("E2").Value = '  -0.48%  '
$ws.Range("D3").Value = '1.893.97'
$ws.Range("E3").Value = '  +1.26%  '
$ws.Range("E4").Value = '  -1.14%  '
$ws.Range("E5").Value = '  -0.44%  '
$ws.Range("E6").Value = '  -1.12%  '
$ws.Range("E7").Value = '  +1.09%  '
$ws.Range("E8").Value = '  -0.43%  '
$ws.Range("E9").Value = '  -0.63%  '
$ws.Range("E10").Value = '  -2.87%  '
$ws.Range("E11").Value = '  -2.15%  '
$ws.Range("E12").Value = '  -1.84%  '
$ws.Range("D13").Value = '1.920.11'
$ws.Range("E13").Value = '  +2.55%  '
$ws.Range("E14").Value = '  -0.41%  '
$ws.Range("E15").Value = '  -0.29%  '
$ws.Range("E16").Value = '  -0.09%  '
$ws.Range("E17").Value = '  -1.17%  '
$ws.Range("E18").Value = '  -0.88%  '
$ws.Range("E19").Value = '  -1.09%  '
$ws.Range("D20").Value = '27.834.57'
$ws.Range("E20").Value = '  +0.01%  '
$ws.Range("E21").Value = '  -2.09%  '
$ws.Range("E22").Value = '  +0.15%  '
$ws.Range("D23").Value = '2.138.50'
$ws.Range("E23").Value = '  +0.37%  '
$ws.Range("E24").Value = '  -1.13%  '
$ws.Range("E25").Value = '  -2.04%  '
$ws.Range("E27").Value = '  -0.95%  '
$ws.Range("E28").Value = '  +4.68%  '
$ws.Range("E29").Value = '  -0.27%  '
$ws.Range("E30").Value = '  -1.74%  '
$ws.Range("E31").Value = '  +0.12%  '
$ws.Range("E32").Value = '  -4.14%  '
$ws.Range("E33").Value = '  -0.70%  '
$ws.Range("E34").Value = '  -0.07%  '
$ws.Range("E35").Value = '  -0.48%  '
$ws.Range("E36").Value = '  -0.72%  '
$ws.Range("E37").Value = '  -7.62%  '
$ws.Range("E38").Value = '  -3.74%  '
$ws.Range("E39").Value = '  -1.75%  '
$ws.Range("E40").Value = '  -2.51%  '
$ws.Range("E41").Value = '  -0.68%  '
$ws.Range("E42").Value = '  -1.82%  '
$ws.Range("E43").Value = '  -0.37%  '
$ws.Range("B44").Value = 'Algorand'
$ws.Range("C44").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("E44").Value = '  -0.83%  '
$ws.Range("B45").Value = 'Quant'
$ws.Range("C45").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("E45").Value = '  +6.63%  '
$ws.Range("E46").Value = '  -0.88%  '
$ws.Range("E47").Value = '  -1.77%  '
$ws.Range("E49").Value = '  -2.15%  '
$ws.Range("E50").Value = '  -0.89%  '
$ws.Range("E51").Value = '  -1.16%  '
